# qumodel trained, some bug fixes in the main function
# Append four new training-set rows to column A, right after the
# existing data (rows 1-139 -> now rows 1-143).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @("дело", "дело", "мис дело", "мис дело")

$startRow = 140
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
